$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.539082
$ws.Range("H2").Value = 58.61724600000001
$ws.Range("I2").Value = 0.224220971665117
$ws.Range("J2").Value = 0.224220971665117
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.151158666666667
$ws.Range("N2").Value = 9.453476
$ws.Range("O2").Value = 0.03114707555614071
$ws.Range("P2").Value = 0.03114707555614071
$ws.Range("Q2").Value = 61.57074758301069
$ws.Range("R2").Value = 554.1367282470961
$ws.Range("S2").Value = 0.006983827545724684
$ws.Range("T2").Value = 0.006983827545724682

$ws.Range("G3").Value = 19.539082
$ws.Range("H3").Value = 58.61724600000001
$ws.Range("I3").Value = 0.224220971665117
$ws.Range("J3").Value = 0.224220971665117
$ws.Range("M3").Value = 5.038243666666667
$ws.Range("O3").Value = 0.04979963650066307
$ws.Range("P3").Value = 0.04979963650066306
$ws.Range("Q3").Value = 98.44265613898069
$ws.Range("R3").Value = 885.9839052508263
$ws.Range("S3").Value = 0.0111661228847483
$ws.Range("T3").Value = 0.0111661228847483

$ws.Range("G4").Value = 19.539082
$ws.Range("H4").Value = 58.61724600000001
$ws.Range("I4").Value = 0.224220971665117
$ws.Range("J4").Value = 0.224220971665117
$ws.Range("M4").Value = 92.91163899999999
$ws.Range("N4").Value = 278.734917
$ws.Range("O4").Value = 0.9183688116343246
$ws.Range("P4").Value = 0.9183688116343246
$ws.Range("Q4").Value = 1815.408133175398
$ws.Range("R4").Value = 16338.67319857858
$ws.Range("S4").Value = 0.205917547291587
$ws.Range("T4").Value = 0.205917547291587

$ws.Range("G5").Value = 19.539082
$ws.Range("H5").Value = 58.61724600000001
$ws.Range("I5").Value = 0.224220971665117
$ws.Range("J5").Value = 0.224220971665117
$ws.Range("M5").Value = 0.06924866666666667
$ws.Range("N5").Value = 0.207746
$ws.Range("O5").Value = 0.0006844763088715736
$ws.Range("P5").Value = 0.0006844763088715734
$ws.Range("Q5").Value = 1.353055376390667
$ws.Range("R5").Value = 12.177498387516
$ws.Range("S5").Value = 0.0001534739430569369
$ws.Range("T5").Value = 0.0001534739430569369

$ws.Range("I6").Value = 0.3010605798326856
$ws.Range("J6").Value = 0.3010605798326856
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.151158666666667
$ws.Range("N6").Value = 9.453476
$ws.Range("O6").Value = 0.03114707555614071
$ws.Range("P6").Value = 0.03114707555614071
$ws.Range("Q6").Value = 82.67079047252622
$ws.Range("R6").Value = 744.0371142527359
$ws.Range("S6").Value = 0.009377156627024192
$ws.Range("T6").Value = 0.009377156627024191

$ws.Range("I7").Value = 0.3010605798326856
$ws.Range("J7").Value = 0.3010605798326856
$ws.Range("M7").Value = 5.038243666666667
$ws.Range("O7").Value = 0.04979963650066307
$ws.Range("P7").Value = 0.04979963650066306
$ws.Range("Q7").Value = 132.1785509953796
$ws.Range("S7").Value = 0.0149927074403466
$ws.Range("T7").Value = 0.0149927074403466

$ws.Range("I8").Value = 0.3010605798326856
$ws.Range("J8").Value = 0.3010605798326856
$ws.Range("M8").Value = 92.91163899999999
$ws.Range("N8").Value = 278.734917
$ws.Range("O8").Value = 0.9183688116343246
$ws.Range("P8").Value = 0.9183688116343246
$ws.Range("Q8").Value = 2437.541061159301
$ws.Range("R8").Value = 21937.86955043371
$ws.Range("S8").Value = 0.2764846469308842
$ws.Range("T8").Value = 0.2764846469308842

$ws.Range("I9").Value = 0.3010605798326856
$ws.Range("J9").Value = 0.3010605798326856
$ws.Range("M9").Value = 0.06924866666666667
$ws.Range("N9").Value = 0.207746
$ws.Range("O9").Value = 0.0006844763088715736
$ws.Range("P9").Value = 0.0006844763088715734
$ws.Range("Q9").Value = 1.816741909272889
$ws.Range("R9").Value = 16.350677183456
$ws.Range("S9").Value = 0.0002060688344306124
$ws.Range("T9").Value = 0.0002060688344306123

$ws.Range("G10").Value = 19.67155566666667
$ws.Range("H10").Value = 59.014667
$ws.Range("I10").Value = 0.2257411748281949
$ws.Range("J10").Value = 0.2257411748281949
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.151158666666667
$ws.Range("N10").Value = 9.453476
$ws.Range("O10").Value = 0.03114707555614071
$ws.Range("P10").Value = 0.03114707555614071
$ws.Range("Q10").Value = 61.98819312583245
$ws.Range("R10").Value = 557.893738132492
$ws.Range("S10").Value = 0.007031177428505757
$ws.Range("T10").Value = 0.007031177428505756

$ws.Range("G11").Value = 19.67155566666667
$ws.Range("H11").Value = 59.014667
$ws.Range("I11").Value = 0.2257411748281949
$ws.Range("J11").Value = 0.2257411748281949
$ws.Range("M11").Value = 5.038243666666667
$ws.Range("O11").Value = 0.04979963650066307
$ws.Range("P11").Value = 0.04979963650066306
$ws.Range("Q11").Value = 99.11009075106412
$ws.Range("R11").Value = 891.9908167595771
$ws.Range("S11").Value = 0.01124182844967674
$ws.Range("T11").Value = 0.01124182844967674

$ws.Range("G12").Value = 19.67155566666667
$ws.Range("H12").Value = 59.014667
$ws.Range("I12").Value = 0.2257411748281949
$ws.Range("J12").Value = 0.2257411748281949
$ws.Range("M12").Value = 92.91163899999999
$ws.Range("N12").Value = 278.734917
$ws.Range("O12").Value = 0.9183688116343246
$ws.Range("P12").Value = 0.9183688116343246
$ws.Range("Q12").Value = 1827.716478669738
$ws.Range("R12").Value = 16449.44830802764
$ws.Range("S12").Value = 0.2073136544639057
$ws.Range("T12").Value = 0.2073136544639057

$ws.Range("G13").Value = 19.67155566666667
$ws.Range("H13").Value = 59.014667
$ws.Range("I13").Value = 0.2257411748281949
$ws.Range("J13").Value = 0.2257411748281949
$ws.Range("M13").Value = 0.06924866666666667
$ws.Range("N13").Value = 0.207746
$ws.Range("O13").Value = 0.0006844763088715736
$ws.Range("P13").Value = 0.0006844763088715734
$ws.Range("Q13").Value = 1.362229001175778
$ws.Range("R13").Value = 12.260061010582
$ws.Range("S13").Value = 0.0001545144861067354
$ws.Range("T13").Value = 0.0001545144861067354

$ws.Range("G14").Value = 21.69639766666667
$ws.Range("H14").Value = 65.08919299999999
$ws.Range("I14").Value = 0.2489772736740025
$ws.Range("J14").Value = 0.2489772736740025
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.151158666666667
$ws.Range("N14").Value = 9.453476
$ws.Range("O14").Value = 0.03114707555614071
$ws.Range("P14").Value = 0.03114707555614071
$ws.Range("Q14").Value = 68.36879154276312
$ws.Range("R14").Value = 615.319123884868
$ws.Range("S14").Value = 0.007754913954886079
$ws.Range("T14").Value = 0.007754913954886077

$ws.Range("G15").Value = 21.69639766666667
$ws.Range("H15").Value = 65.08919299999999
$ws.Range("I15").Value = 0.2489772736740025
$ws.Range("J15").Value = 0.2489772736740025
$ws.Range("M15").Value = 5.038243666666667
$ws.Range("O15").Value = 0.04979963650066307
$ws.Range("P15").Value = 0.04979963650066306
$ws.Range("Q15").Value = 109.3117381335648
$ws.Range("R15").Value = 983.805643202083
$ws.Range("S15").Value = 0.01239897772589143
$ws.Range("T15").Value = 0.01239897772589143

$ws.Range("G16").Value = 21.69639766666667
$ws.Range("H16").Value = 65.08919299999999
$ws.Range("I16").Value = 0.2489772736740025
$ws.Range("J16").Value = 0.2489772736740025
$ws.Range("M16").Value = 92.91163899999999
$ws.Range("N16").Value = 278.734917
$ws.Range("O16").Value = 0.9183688116343246
$ws.Range("P16").Value = 0.9183688116343246
$ws.Range("Q16").Value = 2015.847867605775
$ws.Range("R16").Value = 18142.63080845198
$ws.Range("S16").Value = 0.2286529629479477
$ws.Range("T16").Value = 0.2286529629479477

$ws.Range("G17").Value = 21.69639766666667
$ws.Range("H17").Value = 65.08919299999999
$ws.Range("I17").Value = 0.2489772736740025
$ws.Range("J17").Value = 0.2489772736740025
$ws.Range("M17").Value = 0.06924866666666667
$ws.Range("N17").Value = 0.207746
$ws.Range("O17").Value = 0.0006844763088715736
$ws.Range("P17").Value = 0.0006844763088715734
$ws.Range("Q17").Value = 1.502446609886444
$ws.Range("R17").Value = 13.522019488978
$ws.Range("S17").Value = 0.0001704190452772888
$ws.Range("T17").Value = 0.0001704190452772888
